$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.11449925658460813,
    -0.049877722579360295,
    -0.0089999995842546809,
    0.28399287994709965,
    -0.0059999995973445408,
    -0.0059999995826416352,
    -0.019999999504742405,
    -0.019999999500245558,
    -0.0059999995715473986,
    -0.0059999995672725959,
    0.017669435316552295,
    -0.0059999995651147664,
    -0.0059999995577344478,
    -0.01199999952285502,
    -0.005999999553867319,
    -0.0059999995523045691,
    -0.0059999995502515446,
    -0.0089999995333851501,
    -0.0089999995855469805,
    -0.0089999995818672573,
    -0.0089999995812872768,
    -0.0089999995808653921,
    -0.037435963308036335,
    -0.041999999393167009,
    -0.041999999389831011,
    -0.050802604457160783,
    -0.0059999995808297513,
    -0.0059999995786261806,
    -0.011999999544825002,
    -0.019999999501198573,
    0.045155041610385638,
    -0.020999999496759791,
    -0.0059999995782975546
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
